# Update workplan + report
# Fill in percentages for a few tasks on the "workload" sheet (column F = member 3)
# and leave the SUM/SUMPRODUCT formulas to recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("workload")

$ws.Range("F8").Value = 100
$ws.Range("F9").Value = 100
$ws.Range("F13").Value = 100
$ws.Range("F15").Value = 33.33
$ws.Range("F22").Value = 100

# Update the selected cell to reflect where editing left off.
[void]$ws.Range("F27").Select()

$excel.Calculate()
